$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("0")
$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.3
$ws.Range("D2").Value = 0.3

# Row 3 ("1")
$ws.Range("B3").Value = 0.5
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.5

# Row 4 ("accuracy")
$ws.Range("B4").Value = 0.4166666666666667
$ws.Range("C4").Value = 0.4166666666666667
$ws.Range("D4").Value = 0.4166666666666667
$ws.Range("E4").Value = 0.4166666666666667

# Row 5 ("macro avg")
$ws.Range("B5").Value = 0.4
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 0.4

# Row 6 ("weighted avg")
$ws.Range("B6").Value = 0.4166666666666667
$ws.Range("C6").Value = 0.4166666666666667
$ws.Range("D6").Value = 0.4166666666666667

# Row 7 ("0")
$ws.Range("B7").Value = 0.3636363636363636
$ws.Range("C7").Value = 0.4
$ws.Range("D7").Value = 0.380952380952381

# Row 8 ("1")
$ws.Range("B8").Value = 0.5384615384615384
$ws.Range("C8").Value = 0.5
$ws.Range("D8").Value = 0.5185185185185186

# Row 9 ("accuracy")
$ws.Range("B9").Value = 0.4583333333333333
$ws.Range("C9").Value = 0.4583333333333333
$ws.Range("D9").Value = 0.4583333333333333
$ws.Range("E9").Value = 0.4583333333333333

# Row 10 ("macro avg")
$ws.Range("B10").Value = 0.451048951048951
$ws.Range("C10").Value = 0.45
$ws.Range("D10").Value = 0.4497354497354498

# Row 11 ("weighted avg")
$ws.Range("B11").Value = 0.4656177156177156
$ws.Range("C11").Value = 0.4583333333333333
$ws.Range("D11").Value = 0.4611992945326279

# Row 22 ("0")
$ws.Range("B22").Value = 0.4285714285714285
$ws.Range("C22").Value = 0.6
$ws.Range("D22").Value = 0.5

# Row 23 ("1")
$ws.Range("B23").Value = 0.6
$ws.Range("C23").Value = 0.4285714285714285
$ws.Range("D23").Value = 0.5

# Row 24 ("accuracy")
$ws.Range("B24").Value = 0.5
$ws.Range("C24").Value = 0.5
$ws.Range("D24").Value = 0.5
$ws.Range("E24").Value = 0.5

# Row 25 ("macro avg")
$ws.Range("B25").Value = 0.5142857142857142
$ws.Range("C25").Value = 0.5142857142857142
$ws.Range("D25").Value = 0.5

# Row 26 ("weighted avg")
$ws.Range("B26").Value = 0.5285714285714286
$ws.Range("C26").Value = 0.5
$ws.Range("D26").Value = 0.5
